$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 358 (pushes old rows 358-434 down to 359-435)
$ws.Rows.Item(358).Insert()

$ws.Range("A358").Value = 3
$ws.Range("B358").Value = "Femacal de La Calera"
$ws.Range("C358").Value = "Coquimbo"
$ws.Range("D358").Value = 44907
$ws.Range("E358").Value = 5
$ws.Range("F358").Value = 100112012
$ws.Range("G358").Value = "Espinaca"
$ws.Range("H358").Value = "Sin especificar"
$ws.Range("I358").Value = "Primera"
$ws.Range("J358").Value = 145
$ws.Range("K358").Value = 6000
$ws.Range("L358").Value = 7000
$ws.Range("M358").Value = 6552
$ws.Range("N358").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O358").Value = "Provincia de Quillota"
$ws.Range("P358").Value = 2184
$ws.Range("Q358").Value = 3
$ws.Range("R358").Value = "Hortaliza"

# Insert a second new row at position 434 (after the first shift, pushes former rows
# 434 (old 433) and 435 (old 434) down to 435 and 436)
$ws.Rows.Item(434).Insert()

$ws.Range("A434").Value = 3
$ws.Range("B434").Value = "Femacal de La Calera"
$ws.Range("C434").Value = "Coquimbo"
$ws.Range("D434").Value = 44911
$ws.Range("E434").Value = 5
$ws.Range("F434").Value = 100112012
$ws.Range("G434").Value = "Espinaca"
$ws.Range("H434").Value = "Sin especificar"
$ws.Range("I434").Value = "Primera"
$ws.Range("J434").Value = 130
$ws.Range("K434").Value = 6000
$ws.Range("L434").Value = 6500
$ws.Range("M434").Value = 6231
$ws.Range("N434").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O434").Value = "Provincia de Quillota"
$ws.Range("P434").Value = 2077
$ws.Range("Q434").Value = 3
$ws.Range("R434").Value = "Hortaliza"
